$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the latest IPO (우진엔텍) above the current row 6 (포스뱅크),
# which shifts all the rows below it down by one.
$ws.Rows(6).Insert()

# Fill in the new row 6 with the 우진엔텍 data.
$ws.Range("A6").Value = "우진엔텍"
$ws.Range("B6").Value = "2024.01.08~01.12"
$ws.Range("C6").Value = "4,300~4,900"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 8858
$ws.Range("F6").Value = "케이비증권"

# The table keeps a fixed window of the most recent 20 IPOs, so the oldest
# row (한선엔지니어링), now pushed down to row 22, drops off the bottom.
$ws.Rows(22).Delete()
